$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert a new "Meta description: ..." paragraph right after the
# first (Heading1) paragraph.
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs.First
$null = $p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(2)
$p2.Style = "Normal"

# Borrow the "<w:r/><w:r>text</w:r>" run layout from an existing body
# paragraph (it now sits at index 4, after our insertion) by transplanting
# its FormattedText - this is the only operation that leaves a genuinely
# empty leading run (<w:r/>) in the freshly-created paragraph.
$sourceForEmptyRun = $d.Paragraphs.Item(4)
$ft = $sourceForEmptyRun.Range.FormattedText
$p2.Range.FormattedText = $ft

# Split off a dedicated run for the first 16 characters ("Meta description")
# by toggling Bold - this creates a real run boundary *without* touching /
# merging the leading empty run, unlike a straight text replace would.
$p2 = $d.Paragraphs.Item(2)
$boldRange = $d.Range($p2.Range.Start, $p2.Range.Start + 16)
$boldRange.Bold = 1

# Now swap the placeholder text for the real copy, one run at a time, so the
# leading empty run is never part of the edited range.
$p2 = $d.Paragraphs.Item(2)
$boldRange = $d.Range($p2.Range.Start, $p2.Range.Start + 16)
$boldRange.Find.Execute($boldRange.Text, $true, $false, $false, $false, $false, `
    $true, 1, $false, "Meta description", 2)

$p2 = $d.Paragraphs.Item(2)
$restStart = $p2.Range.Start + 16
$restEnd = $p2.Range.End - 1
$restRange = $d.Range($restStart, $restEnd)
$restRange.Find.Execute($restRange.Text, $true, $false, $false, $false, $false, `
    $true, 1, $false, ": Read our review of Da Vinci Extreme, a 5-reel slot machine game with tumbling reels and a free spins bonus, and play for free.", 2)

# ---------------------------------------------------------------------------
# Change 2: at the very end of the document, drop the duplicated bold
# "Play Da Vinci Extreme..." paragraph and rewrite the italic paragraph's
# text into the DALL-E image prompt.
# ---------------------------------------------------------------------------

$count = $d.Paragraphs.Count
$pBold = $d.Paragraphs.Item($count - 1)
$pBold.Range.Delete()

$count = $d.Paragraphs.Count
$pItalic = $d.Paragraphs.Item($count)
$oldPrompt = "Read our review of Da Vinci Extreme, a 5-reel slot machine game with tumbling reels and a free spins bonus, and play for free."
$newPrompt = "Prompt: Create a feature image for Da Vinci Extreme that is in a cartoon style and features a happy Maya warrior with glasses. For this feature image, DALLE could create a cartoon-style image of a happy Maya warrior with glasses standing next to the game's logo or a slot machine. The Maya warrior could be adorned with traditional clothing and accessories, such as a headdress, necklace, and bracelets. They could be holding a tablet or smartphone, indicating that the game can be played on mobile devices. The background could be a Renaissance-inspired painting or artwork, tying in with the game's theme. The overall image should be bright, colorful, and eye-catching to appeal to players of all ages and genders."
$pItalic.Range.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
    $true, 1, $false, $newPrompt, 2)

Write-Output "done"
